$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.293.89'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").Value = '2.612.06'
$ws.Range("E3").Value = '  +7.99%  '

$ws.Range("E4").Value = '  -0.13%  '

$cell = $ws.Range("D5")
$cell.Formula = "'313.99"
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.49%  '

$cell = $ws.Range("D6")
$cell.Formula = "'101.06"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +4.28%  '

$cell = $ws.Range("D7")
$cell.Formula = "'0.601"
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +6.26%  '

$ws.Range("E8").Value = '  +0.00%  '

$cell = $ws.Range("D9")
$cell.Formula = "'0.584"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +14.00%  '

$cell = $ws.Range("D10")
$cell.Formula = "'38.85"
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +12.56%  '

$cell = $ws.Range("D11")
$cell.Formula = "'54.34"
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '

$ws.Range("E12").Value = '  +6.47%  '

$ws.Range("E13").Value = '  +16.02%  '

$ws.Range("D14").Value = '3.012.05'
$ws.Range("E14").Value = '  +7.97%  '

$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("D16").Value = '2.606.49'
$ws.Range("E16").Value = '  +7.17%  '

$cell = $ws.Range("D17")
$cell.Formula = "'0.911"
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +8.25%  '

$cell = $ws.Range("D18")
$cell.Formula = "'15.12"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +5.54%  '

$ws.Range("D19").Value = '46.526.50'
$ws.Range("E19").Value = '  +1.97%  '

$cell = $ws.Range("D20")
$cell.Formula = "'13.32"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("E21").Value = '  +7.79%  '

$ws.Range("E22").Value = '  +9.70%  '

$cell = $ws.Range("D23")
$cell.Formula = "'71.02"
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +5.01%  '

$cell = $ws.Range("D24")
$cell.Formula = "'255.28"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +4.54%  '

$cell = $ws.Range("D25")
$cell.Formula = "'3.09"
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +10.16%  '

$ws.Range("E26").Value = '  +14.21%  '

$cell = $ws.Range("D27")
$cell.Formula = "'28.02"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +30.91%  '

$cell = $ws.Range("D28")
$cell.Formula = "'0.999"
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '

$cell = $ws.Range("D29")
$cell.Formula = "'10.62"
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +8.11%  '

$cell = $ws.Range("D30")
$cell.Formula = "'40.91"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.10%  '

$ws.Range("E31").Value = '  +2.51%  '

$cell = $ws.Range("D32")
$cell.Formula = "'6.19"
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +10.85%  '

$cell = $ws.Range("D33")
$cell.Formula = "'3.73"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.66%  '

$ws.Range("E34").Value = '  +14.30%  '

$cell = $ws.Range("D35")
$cell.Formula = "'2.88"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +4.91%  '

$cell = $ws.Range("D36")
$cell.Formula = "'154.07"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +4.13%  '

$ws.Range("E37").Value = '  +7.87%  '

$cell = $ws.Range("D38")
$cell.Formula = "'0.119"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +5.26%  '

$ws.Range("E39").Value = '  +5.71%  '

$cell = $ws.Range("D40")
$cell.Formula = "'17.07"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +12.10%  '

$cell = $ws.Range("D41")
$cell.Formula = "'4.25"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +9.59%  '

$cell = $ws.Range("D42")
$cell.Formula = "'3.66"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +11.41%  '

$ws.Range("E43").Value = '  +8.67%  '

$cell = $ws.Range("D44")
$cell.Formula = "'21.31"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +43.55%  '

$ws.Range("D45").Value = '2.032.89'
$ws.Range("E45").Value = '  +3.59%  '

$ws.Range("E46").Value = '  +0.00%  '

$cell = $ws.Range("D47")
$cell.Formula = "'91.33"
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.53%  '

$cell = $ws.Range("D48")
$cell.Formula = "'111.99"
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +11.56%  '

$ws.Range("E50").Value = '  +2.31%  '

$cell = $ws.Range("D51")
$cell.Formula = "'0.203"
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +9.25%  '
